$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row with t+3 / t+4 prediction columns, matching formatting of existing header cells
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:E1").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4

# Update predicted-variable values for columns B:E across all data rows
$data = @"
B2=-0.3262056301079485
C2=-0.3112980140467429
D2=-0.2968247675687409
E2=-0.2831294276583772
B3=0.1973430043990916
C3=0.2065416847512486
D3=0.2166372530487319
E3=0.2267245108538933
B4=0.1302510715930824
C4=0.1418205654855191
D4=0.1543698210598723
E4=0.1669317699465675
B5=-0.08662491400362807
C5=-0.07704090021470952
D5=-0.06782744480529213
E5=-0.05947266192440955
B6=0.1267911413067892
C6=0.1425763902569039
D6=0.1577757995213483
E6=0.1716403738725467
B7=-0.4622718219935082
C7=-0.4527558810974688
D7=-0.4427275793353809
E7=-0.4326515851592884
B8=-0.2721202125849562
C8=-0.2603289572191466
D8=-0.2477026304728691
E8=-0.2349847685432591
B9=-0.4215539029919932
C9=-0.4114738013147691
D9=-0.3991547334406629
E9=-0.385583653891701
B10=0.3433447984396565
C10=0.3580686953613876
D10=0.3719365335256363
E10=0.3843220301528791
B11=-0.2487543778038808
C11=-0.2399055122828032
D11=-0.230912730326157
E11=-0.2222892167707231
B12=-0.1097844704537383
C12=-0.0873066370082155
D12=-0.06702825810620978
E12=-0.04937651895594366
B13=-0.01171595900957665
C13=-0.01191385433495715
D13=-0.01156541023472136
E13=-0.0112101582789857
B14=0.01149504879265671
C14=0.05155802922678206
D14=0.08801552103037588
E14=0.1200856704240608
B15=-0.05819741262321203
C15=-0.01369572637148195
D15=0.02558501525861567
E15=0.05909453672944824
B16=0.2858159185328134
C16=0.3472103327010492
D16=0.4022562303047558
E16=0.4500718595671089
B17=0.563124343860992
C17=0.598514994193861
D17=0.6285717138606319
E17=0.652818407669788
B18=0.05151022774580531
C18=0.03878916555125164
D18=0.02770627094762623
E18=0.01770400497480018
B19=0.3801050985279987
C19=0.3958043559957307
D19=0.409721124107925
E19=0.4211276169770829
B20=0.2354360921909205
C20=0.2934706472736356
D20=0.3446824427181091
E20=0.3883515187112792
B21=0.4604852122243949
C21=0.517878593557578
D21=0.5671489504056036
E21=0.608052850793377
B22=0.3407563088633724
C22=0.377430769596674
D22=0.4090313796388642
E22=0.4350154813545994
B23=-0.09244117645539568
C23=-0.0608179811506077
D23=-0.03397398560891771
E23=-0.01216876333360331
B24=4.537685292080264
C24=4.579137554806769
D24=4.537864491830653
E24=4.429727826019735
B25=0.5276714520015621
C25=0.488382454367825
D25=0.4551840160358246
E25=0.4255583837517669
B26=0.3865084971551981
C26=0.3615177383265028
D26=0.3398179294592206
E26=0.3178135699851935
B27=0.3347112899583037
C27=0.303409187322828
D27=0.2758292037665361
E27=0.248481584850865
B28=1.101682590302006
C28=1.069211801959468
D28=1.037288312556877
E28=1.005714189489717
B29=5.880675973654099
C29=5.465989120357457
D29=5.057806472844939
E29=4.661229097481694
B30=1.026257237776456
C30=0.9723697196820085
D30=0.9242096999879832
E30=0.879851145615141
B31=-0.132021831180891
C31=-0.1983294097244325
D31=-0.2546211676007799
E31=-0.303876126658982
B32=0.8160840418720325
C32=0.7801612285857227
D32=0.748210789192325
E32=0.717327399663262
B33=0.943321890479472
C33=0.912811341700061
D33=0.8868777381331967
E33=0.8628002275848395
B34=-0.6225044766211654
C34=-0.6563600928369668
D34=-0.685891609971865
E34=-0.7133167669137845
B35=0.8227091330419488
C35=0.8136177197400545
D35=0.8057086695402788
E35=0.7983758155621118
B36=0.785324575538821
C36=0.7676957230286918
D36=0.7527353103908232
E36=0.739766476328692
B37=0.7616482074449005
C37=0.7417196672734563
D37=0.7246931780510263
E37=0.7099087826312352
B38=0.7368169998210876
C38=0.7176524242778116
D38=0.700922380920297
E38=0.6859547615923561
B39=0.5852317773954258
C39=0.580898655261708
D39=0.5776250229312484
E39=0.5747992602722133
B40=0.7552084469826726
C40=0.7521349625251679
D40=0.7494396738585015
E40=0.7465407144997647
B41=0.5626100458938021
C41=0.5543715805352603
D41=0.5483935439029464
E41=0.5438690680778824
B42=0.7179076915742605
C42=0.6904290198476798
D42=0.6662060303902337
E42=0.6446789416152406
B43=0.7241739736172695
C43=0.7091434782750188
D43=0.6964702833621887
E43=0.6854170303110387
B44=0.6829084950922456
C44=0.674411792363956
D44=0.6678673715680018
E44=0.6625382763182925
B45=0.6770882469038987
C45=0.6580905976512383
D45=0.6422992906277425
E45=0.6290737183597467
B46=-1.259294341841793
C46=-1.263811510086398
D46=-1.267708564968823
E46=-1.270867090200296
B47=-0.9741234816865039
C47=-0.9795097922288276
D47=-0.9841360423357495
E47=-0.9879821072228756
B48=-0.8669027100444782
C48=-0.8735394605770974
D48=-0.878526254468243
E48=-0.8820532658286647
B49=-0.6387071340702309
C49=-0.6425263145464726
D49=-0.6451204534471169
E49=-0.6467585949912471
B50=-0.04734612250935647
C50=-0.0488315234424316
D50=-0.04986712173913663
E50=-0.05084090092411544
B51=-0.8597814277296713
C51=-0.8654141101463644
D51=-0.8697001213243605
E51=-0.8727522574982215
B52=-0.8597814277296713
C52=-0.8654141101463644
D52=-0.8697001213243605
E52=-0.8727522574982215
B53=-1.083552313285897
C53=-1.097302253505303
D53=-1.108638927503274
E53=-1.117691261675691
B54=-0.1874860143018053
C54=-0.1871342512186376
D54=-0.1861820668096902
E54=-0.1850206658405746
B55=-0.9947846777069013
C55=-0.9993792754448271
D55=-1.003344435830261
E55=-1.006677793986441
B56=-0.8878718671948455
C56=-0.8844730449215408
D56=-0.8820658690662952
E56=-0.8805137556689066
B57=-0.9468508169596056
C57=-0.9341162932906286
D57=-0.9228142763993664
E57=-0.9130028463906261
B58=-1.15126224138174
C58=-1.128542268251601
D58=-1.108593711885405
E58=-1.091286089897769
B59=-0.8600764908849532
C59=-0.8483927370979433
D59=-0.8373464950248912
E59=-0.8270450846581709
B60=-0.5183189354500672
C60=-0.4981670011442724
D60=-0.4797639111915284
E60=-0.4633679861010591
B61=0.3693746975788355
C61=0.3713187932361238
D61=0.3738096984249569
E61=0.3761934843145934
B62=-1.229297507291419
C62=-1.209750680733107
D62=-1.192594656859132
E62=-1.177787945402543
B63=-0.7657356006018912
C63=-0.7400286323080764
D63=-0.7155447674589617
E63=-0.6925436139656344
B64=-0.9088214832527715
C64=-0.9027772512895821
D64=-0.8956121157317797
E64=-0.8877465956533145
B65=-0.1312481098718712
C65=-0.1119708445767711
D65=-0.09441090829008827
E65=-0.07892842271157849
B66=-0.8149438709111899
C66=-0.7911518241820487
D66=-0.7708609654517908
E66=-0.7540318667592777
B67=-0.7951699883843418
C67=-0.7628194820164587
D67=-0.7361516981732094
E67=-0.7148266051416914
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split "=", 2
    $addr = $parts[0]
    $val = [double]$parts[1]
    $ws.Range($addr).Value = $val
}
